# Code fix for Org
$wb = $excel.ActiveWorkbook

# --- Portfolio sheet: fill in row 2 data ---
$wsPortfolio = $wb.Worksheets.Item("Portfolio")
$rowRange = $wsPortfolio.Range("A2:N2")

# Force text formatting first so numeric-looking values (15, 10, 26556,
# 7362255383) are stored as shared-string text rather than numbers.
$rowRange.NumberFormat = "@"

# Values are written in the same order the source automation run entered
# them (location/owner fields first, then the account/org/address fields)
# so new shared-string entries land in the same order as the target file.
$wsPortfolio.Range("D2").Value = "15"
$wsPortfolio.Range("E2").Value = "10"
$wsPortfolio.Range("L2").Value = "Joey"
$wsPortfolio.Range("M2").Value = "jayne.keebler@hotmail.com"
$wsPortfolio.Range("N2").Value = "7362255383"
$wsPortfolio.Range("A2").Value = "Automation portfolio 4377747"
$wsPortfolio.Range("B2").Value = "Government"
$wsPortfolio.Range("F2").Value = "26556"
$wsPortfolio.Range("H2").Value = "Ohio"
$wsPortfolio.Range("I2").Value = "88741 Lucas Locks"
$wsPortfolio.Range("J2").Value = "Durganberg"
$wsPortfolio.Range("K2").Value = "59989-8155"
$wsPortfolio.Range("C2").Value = "100 RESILIENT CITIES"
$wsPortfolio.Range("G2").Value = "United States"

# Drop the temporary text number-format back to the default style so the
# cells don't carry a leftover style index.
$rowRange.Style = "Normal"

$wsPortfolio.Range("C1").Select()

# --- Faculty sheet: add "Org" header column ---
$wsFaculty = $wb.Worksheets.Item("Faculty")
$wsFaculty.Range("I1").Value = "Org"
$wsFaculty.Range("A4").Select()

# --- Hsr sheet: add "Org" header column, becomes active/selected ---
$wsHsr = $wb.Worksheets.Item("Hsr")
$wsHsr.Range("F1").Value = "Org"
$wsHsr.Activate()
$wsHsr.Range("F1").Select()
